$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for column C (nombre_aides) and D (montant_total)
# for the 2020-09-04 data refresh.
$updates = @(
    @{Row=2; C=39031; D=56426421},
    @{Row=3; C=93550; D=137119130},
    @{Row=4; C=31924; D=47274264},
    @{Row=5; C=8988; D=13358746},
    @{Row=6; C=2096; D=3115971},
    @{Row=7; C=175; D=257593},
    @{Row=12; C=42417; D=57496974},
    @{Row=13; C=9954; D=14390743},
    @{Row=14; C=26577; D=38958476},
    @{Row=15; C=8490; D=12600478},
    @{Row=16; C=2216; D=3292539},
    @{Row=17; C=433; D=638623},
    @{Row=20; C=10447; D=13818095},
    @{Row=21; C=13751; D=19843341},
    @{Row=22; C=32382; D=47508431},
    @{Row=23; C=10460; D=15546408},
    @{Row=24; C=2712; D=4032771},
    @{Row=25; C=546; D=813092},
    @{Row=26; C=37; D=54953},
    @{Row=27; C=11965; D=15968134},
    @{Row=28; C=7907; D=11441882},
    @{Row=29; C=23098; D=33906358},
    @{Row=30; C=7977; D=11864392},
    @{Row=31; C=2020; D=3014251},
    @{Row=32; C=381; D=568915},
    @{Row=34; C=8524; D=11257893},
    @{Row=35; C=3377; D=4877500},
    @{Row=36; C=8063; D=11776926},
    @{Row=37; C=3240; D=4802961},
    @{Row=41; C=2546; D=3440228},
    @{Row=42; C=17829; D=25778696},
    @{Row=43; C=52449; D=76868936},
    @{Row=44; C=19373; D=28769257},
    @{Row=45; C=5761; D=8575187},
    @{Row=46; C=1267; D=1891045},
    @{Row=48; C=11; D=16500},
    @{Row=50; C=17202; D=22849422},
    @{Row=51; C=2166; D=3144228},
    @{Row=52; C=7301; D=10729727},
    @{Row=53; C=2447; D=3654572},
    @{Row=54; C=777; D=1160415},
    @{Row=57; C=7400; D=10179114},
    @{Row=58; C=1214; D=2106300},
    @{Row=59; C=2949; D=5077033},
    @{Row=60; C=1163; D=2013106},
    @{Row=61; C=403; D=704383},
    @{Row=63; C=31; D=64500},
    @{Row=64; C=1781; D=2848310},
    @{Row=65; C=15909; D=22975720},
    @{Row=66; C=45960; D=67238303},
    @{Row=67; C=16055; D=23854044},
    @{Row=68; C=4677; D=6965288},
    @{Row=69; C=974; D=1447802},
    @{Row=71; C=15; D=21287},
    @{Row=73; C=15477; D=20380459},
    @{Row=74; C=55265; D=80409209},
    @{Row=75; C=153969; D=226789431},
    @{Row=76; C=66221; D=98668728},
    @{Row=77; C=21244; D=31743884},
    @{Row=78; C=5113; D=7637903},
    @{Row=79; C=291; D=431670},
    @{Row=85; C=54151; D=73541913},
    @{Row=86; C=4797; D=6949364},
    @{Row=87; C=11945; D=17544645},
    @{Row=88; C=3978; D=5926958},
    @{Row=89; C=1378; D=2058789},
    @{Row=90; C=298; D=444512},
    @{Row=93; C=5587; D=7505560},
    @{Row=94; C=1673; D=2411699},
    @{Row=95; C=5377; D=7922340},
    @{Row=96; C=1995; D=2970426},
    @{Row=101; C=3716; D=4923558},
    @{Row=102; C=748; D=1286775},
    @{Row=103; C=465; D=831027},
    @{Row=106; C=31; D=61500},
    @{Row=107; C=11143; D=16158284},
    @{Row=108; C=29874; D=43869885},
    @{Row=109; C=10015; D=14888226},
    @{Row=110; C=2765; D=4122580},
    @{Row=112; C=58; D=87000},
    @{Row=114; C=10044; D=13255019},
    @{Row=115; C=31405; D=45274835},
    @{Row=116; C=67787; D=99180797},
    @{Row=117; C=21824; D=32427399},
    @{Row=118; C=6199; D=9234341},
    @{Row=119; C=1166; D=1742600},
    @{Row=124; C=26438; D=35281282},
    @{Row=125; C=37241; D=53734144},
    @{Row=126; C=78935; D=115409102},
    @{Row=127; C=24405; D=36222952},
    @{Row=128; C=6552; D=9736623},
    @{Row=129; C=1312; D=1951311},
    @{Row=133; C=32601; D=43267777},
    @{Row=134; C=13707; D=19842737},
    @{Row=135; C=33132; D=48655323},
    @{Row=136; C=11727; D=17423816},
    @{Row=137; C=3042; D=4534241},
    @{Row=138; C=516; D=767990},
    @{Row=141; C=11075; D=14759097},
    @{Row=142; C=36389; D=52554984},
    @{Row=143; C=83887; D=122890569},
    @{Row=144; C=24998; D=37137458},
    @{Row=145; C=6557; D=9783996},
    @{Row=146; C=1498; D=2229230},
    @{Row=149; C=30023; D=40470887}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
